$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4220.727
$ws.Range("I86").Value = 2828.5715
$ws.Range("J86").Value = 4595.5386
$ws.Range("K86").Value = 2828.5715
$ws.Range("L86").Value = 4595.5386
$ws.Range("M86").Value = -1705.5715
$ws.Range("N86").Value = -6841.5386

# Sheet ALC, Row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 32350
$ws.Range("J87").Value = 32350
$ws.Range("L87").Value = 32350
$ws.Range("N87").Value = -34846

# Sheet ALC, Row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4220.727
$ws.Range("I89").Value = 2828.5715
$ws.Range("J89").Value = 4595.5386
$ws.Range("K89").Value = 14142.8575
$ws.Range("L89").Value = 22977.693
$ws.Range("M89").Value = -8526.8575
$ws.Range("N89").Value = -34209.693

# Sheet ALC, Row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 32350
$ws.Range("J90").Value = 32350
$ws.Range("L90").Value = 97050
$ws.Range("N90").Value = -109530

# Sheet ALC, Row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3121
$ws.Range("I98").Value = 793.7692
$ws.Range("J98").Value = 8163.3335
$ws.Range("K98").Value = 793.7692
$ws.Range("L98").Value = 8163.3335
$ws.Range("M98").Value = 704.2308
$ws.Range("N98").Value = -11159.3335

# Sheet ALC, Row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3121
$ws.Range("I122").Value = 793.7692
$ws.Range("J122").Value = 8163.3335
$ws.Range("K122").Value = 2381.3076
$ws.Range("L122").Value = 24490.0005
$ws.Range("M122").Value = 68.69239999999991
$ws.Range("N122").Value = -29390.0005

# Sheet ALC, Row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1448.7858
$ws.Range("I127").Value = 357
$ws.Range("J127").Value = 2055.3333
$ws.Range("K127").Value = 1071
$ws.Range("L127").Value = 6165.999899999999
$ws.Range("M127").Value = 3889
$ws.Range("N127").Value = -16085.9999

# Sheet ALC, Row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 979.2442
$ws.Range("J129").Value = 1022.58026
$ws.Range("L129").Value = 3067.74078
$ws.Range("N129").Value = -13067.74078

# Sheet ARM, Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1577.7778
$ws.Range("I45").Value = 1577.7778
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1577.7778
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1200.7778
$ws.Range("N45").ClearContents()

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8107.242
$ws.Range("I132").Value = 7037.5
$ws.Range("K132").Value = 21112.5
$ws.Range("M132").Value = -18582.5

# Sheet BSM, Row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2775
$ws.Range("I94").Value = 1330
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 1330
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -879
$ws.Range("N94").Value = -10902

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11670.2
$ws.Range("I132").Value = 9353.833000000001
$ws.Range("J132").Value = 15144.75
$ws.Range("K132").Value = 28061.499
$ws.Range("L132").Value = 45434.25
$ws.Range("M132").Value = -25531.499
$ws.Range("N132").Value = -50494.25

# Sheet CRP, Row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8143.161
$ws.Range("I134").Value = 8776.714
$ws.Range("J134").Value = 7621.4116
$ws.Range("K134").Value = 26330.142
$ws.Range("L134").Value = 22864.2348
$ws.Range("M134").Value = -23795.142
$ws.Range("N134").Value = -27934.2348

# Sheet CUL, Row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3060.5715
$ws.Range("I34").Value = 444
$ws.Range("J34").Value = 3496.6667
$ws.Range("K34").Value = 1332
$ws.Range("L34").Value = 10490.0001
$ws.Range("M34").Value = -1248
$ws.Range("N34").Value = -10658.0001

# Sheet CUL, Row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2884.6155
$ws.Range("I55").Value = 750
$ws.Range("J55").Value = 3272.7273
$ws.Range("K55").Value = 2250
$ws.Range("L55").Value = 9818.1819
$ws.Range("M55").Value = -2073
$ws.Range("N55").Value = -10172.1819

# Sheet CUL, Row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2303.3333
$ws.Range("I136").Value = 1575
$ws.Range("J136").Value = 3760
$ws.Range("K136").Value = 4725
$ws.Range("L136").Value = 11280
$ws.Range("M136").Value = 375
$ws.Range("N136").Value = -21480

# Sheet GSM, Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1280.1538
$ws.Range("I102").Value = 1197.5264
$ws.Range("J102").Value = 1504.4286
$ws.Range("K102").Value = 1197.5264
$ws.Range("L102").Value = 1504.4286
$ws.Range("M102").Value = 424.4736
$ws.Range("N102").Value = -4748.4286

# Sheet GSM, Row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1900
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -3230

# Sheet LTW, Row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3436.6667
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3436.6667
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3436.6667
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -3660.6667

# Sheet LTW, Row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5537.375
$ws.Range("I40").Value = 4779.8
$ws.Range("J40").Value = 6800
$ws.Range("K40").Value = 4779.8
$ws.Range("L40").Value = 6800
$ws.Range("M40").Value = -4643.8
$ws.Range("N40").Value = -7072

# Sheet LTW, Row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

# Sheet LTW, Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4890.6665
$ws.Range("I122").Value = 5089.643
$ws.Range("J122").Value = 2105
$ws.Range("K122").Value = 15268.929
$ws.Range("L122").Value = 6315
$ws.Range("M122").Value = -12818.929
$ws.Range("N122").Value = -11215

# Sheet LTW, Row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3436.6667
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3436.6667
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10310.0001
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -15250.0001

# Sheet LTW, Row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5725.722
$ws.Range("I132").Value = 6296.846
$ws.Range("K132").Value = 18890.538
$ws.Range("M132").Value = -16360.538

# Sheet WVR, Row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 590620.8
$ws.Range("I122").Value = 1430435.6
$ws.Range("J122").Value = 2750.5
$ws.Range("K122").Value = 4291306.800000001
$ws.Range("L122").Value = 8251.5
$ws.Range("M122").Value = -4288856.800000001
$ws.Range("N122").Value = -13151.5

# Sheet WVR, Row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 313327.72
$ws.Range("I126").Value = 417182.38
$ws.Range("J126").Value = 1763.75
$ws.Range("K126").Value = 1251547.14
$ws.Range("L126").Value = 5291.25
$ws.Range("M126").Value = -1249077.14
$ws.Range("N126").Value = -10231.25

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15243.25
$ws.Range("I132").Value = 18658
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 55974
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -53444
$ws.Range("N132").Value = -20057
